# Update countries & provincias Spain
# Applies the latest COVID country-data update:
#  - refresh the case counters for the countries whose figures moved
#  - re-sort the table (rows 4:219) descending by "Casos totales" (col B),
#    which is how the sheet is always kept ordered
#  - bump the "Datos actualizados a ..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# country name -> new values for columns B,C,D,E,F,G,H
# (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$updates = @{
    "Estados Unidos"   = @(6824441, 36294, 4101452, 2521748, 0, 1044, 201241)
    "Brasil"           = @(4419083, 34784, 3720312, 564665, 0, 899, 134106)
    "Colombia"         = @(736377, 7787, 610078, 102821, 0, 190, 23478)
    "Canada"           = @(139625, 822, 122175, 8257, 0, 5, 9193)
    "Egipto"           = @(101500, 160, 86549, 9255, 0, 17, 5696)
    "Guatemala"        = @(82924, 240, 72562, 7353, 0, 25, 3009)
    "Japon"            = @(76448, 490, 68532, 6455, 0, 10, 1461)
    "Barein"           = @(62484, 841, 55444, 6824, 0, 3, 216)
    "Nigeria"          = @(56604, 126, 47872, 7641, 0, 3, 1091)
    "Camerun"          = @(20303, 32, 18837, 1051, 0, 0, 415)
    "Bulgaria"         = @(18390, 174, 13241, 4410, 0, 3, 739)
    "Noruega"          = @(12498, 105, 10371, 1862, 0, 0, 265)
    "Guinea"           = @(10154, 43, 9612, 479, 0, 0, 63)
    "Birmania"         = @(3821, 319, 908, 2873, 0, 5, 40)
    "Trinidad yTobago" = @(3327, 104, 810, 2459, 0, 2, 58)
    "Togo"             = @(1608, 13, 1230, 338, 0, 0, 40)
    "Curazao"          = @(192, 23, 68, 123, 0, 0, 1)
}

$searchRange = $ws.Range("A4:A219")
$xlWhole = 1
$xlValues = -4163

foreach ($name in $updates.Keys) {
    $vals = $updates[$name]
    $found = $searchRange.Find($name, $null, $xlValues, $xlWhole)
    if ($found -eq $null) {
        Write-Host ("WARNING: country not found: " + $name)
        continue
    }
    $row = $found.Row
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# Re-sort the data block (rows 4..219) descending by column B (Casos totales),
# exactly as the source feed keeps it ranked.
$dataRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$dataRange.Sort($keyRange, 2, $null, $null, $null, $null, $null, 1)

# Update the "last updated" banner.
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 00:51"

$wb.Save()
